# Adds a new "localdb" command-category column to the hidden '#system'
# worksheet (used to back the drop-down/autocomplete defined-names), and
# registers the "localdb" entry in the "target" (category list) column.
#
# Net effect mirrors what the Excel UI does when a user inserts a new
# column at N (shifting macro..xml from N..AC to O..AD), fills the new
# column with the localdb function names, inserts "localdb" into the
# alphabetically-sorted target list in column A (pushing macro..xml down
# one row), and fixes up all the defined names that pointed at the
# shifted columns/ranges - plus adds the new "localdb" defined name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Insert a new column before column N (14). This shifts the existing
#    macro..xml columns (N..AC) one column to the right (O..AD), exactly
#    like Excel's "Insert Column" command; it does not touch column A.
# ---------------------------------------------------------------------
$ws.Columns.Item(14).Insert()

# ---------------------------------------------------------------------
# 2) Populate the newly-inserted column N with the "localdb" category.
# ---------------------------------------------------------------------
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------
# 3) Insert "localdb" into the alphabetically-sorted "target" list that
#    lives in column A (A2:A29 -> A2:A30), shifting macro..xml down by
#    one row. Only column A moves - the other category columns keep
#    their own row numbers, so this is done with manual cell copies
#    rather than a real row insert (which would shift every column).
# ---------------------------------------------------------------------
for ($r = 29; $r -ge 14; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value2
}
$ws.Cells.Item(14, 1).Value = "localdb"

# ---------------------------------------------------------------------
# 4) Fix up the workbook-level defined names so they keep pointing at
#    the right ranges after the column insert / target-list growth, and
#    add the new "localdb" name.
# ---------------------------------------------------------------------
$wb.Names.Item("mail").RefersTo = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("macro").RefersTo = "='#system'!`$O`$2:`$O`$4"

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
